$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (tab name) to reflect new date
$ws.Name = "Through 2022-03-07"

# Update the "March (through 03-06)" label cell
$ws.Range("A4").Value = "March (through 03-07)"

# Update March row (row 4) values
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 18
$ws.Range("I4").Value = 36

# Update Total row (row 5) values
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 97
$ws.Range("D5").Value = 141
$ws.Range("E5").Value = 149
$ws.Range("F5").Value = 86
$ws.Range("G5").Value = 157
$ws.Range("H5").Value = 360
$ws.Range("I5").Value = 337
